$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.636579666666667
$ws.Range("H2").Value = 19.909739
$ws.Range("I2").Value = 0.1201574291771603
$ws.Range("J2").Value = 0.1201574291771603
$ws.Range("M2").Value = 38.45264233333334
$ws.Range("N2").Value = 115.357927
$ws.Range("O2").Value = 0.2975040117664333
$ws.Range("P2").Value = 0.2975040117664332
$ws.Range("Q2").Value = 255.1940242390059
$ws.Range("R2").Value = 2296.746218151053
$ws.Range("S2").Value = 0.03574731722374627
$ws.Range("T2").Value = 0.03574731722374627
$ws.Range("G3").Value = 6.636579666666667
$ws.Range("H3").Value = 19.909739
$ws.Range("I3").Value = 0.1201574291771603
$ws.Range("J3").Value = 0.1201574291771603
$ws.Range("O3").Value = 0.3694391181876273
$ws.Range("P3").Value = 0.3694391181876272
$ws.Range("Q3").Value = 316.8987696059284
$ws.Range("R3").Value = 2852.088926453355
$ws.Range("S3").Value = 0.04439085467890238
$ws.Range("T3").Value = 0.04439085467890238
$ws.Range("G4").Value = 6.636579666666667
$ws.Range("H4").Value = 19.909739
$ws.Range("I4").Value = 0.1201574291771603
$ws.Range("J4").Value = 0.1201574291771603
$ws.Range("M4").Value = 18.63107466666667
$ws.Range("N4").Value = 55.893224
$ws.Range("O4").Value = 0.1441466469015163
$ws.Range("P4").Value = 0.1441466469015162
$ws.Range("Q4").Value = 123.6466113009485
$ws.Range("R4").Value = 1112.819501708536
$ws.Range("S4").Value = 0.01732029051619408
$ws.Range("T4").Value = 0.01732029051619408
$ws.Range("G5").Value = 6.636579666666667
$ws.Range("H5").Value = 19.909739
$ws.Range("I5").Value = 0.1201574291771603
$ws.Range("J5").Value = 0.1201574291771603
$ws.Range("M5").Value = 24.41680433333333
$ws.Range("N5").Value = 73.25041299999999
$ws.Range("O5").Value = 0.1889102231444233
$ws.Range("P5").Value = 0.1889102231444233
$ws.Range("Q5").Value = 162.0440671635785
$ws.Range("R5").Value = 1458.396604472207
$ws.Range("S5").Value = 0.02269896675831759
$ws.Range("T5").Value = 0.02269896675831759
$ws.Range("I6").Value = 0.00477103065019021
$ws.Range("J6").Value = 0.00477103065019021
$ws.Range("M6").Value = 38.45264233333334
$ws.Range("N6").Value = 115.357927
$ws.Range("O6").Value = 0.2975040117664333
$ws.Range("P6").Value = 0.2975040117664332
$ws.Range("Q6").Value = 10.13286086201578
$ws.Range("R6").Value = 91.19574775814201
$ws.Range("S6").Value = 0.001419400758692202
$ws.Range("T6").Value = 0.001419400758692202
$ws.Range("I7").Value = 0.00477103065019021
$ws.Range("J7").Value = 0.00477103065019021
$ws.Range("O7").Value = 0.3694391181876273
$ws.Range("P7").Value = 0.3694391181876272
$ws.Range("S7").Value = 0.001762605356252413
$ws.Range("T7").Value = 0.001762605356252413
$ws.Range("I8").Value = 0.00477103065019021
$ws.Range("J8").Value = 0.00477103065019021
$ws.Range("M8").Value = 18.63107466666667
$ws.Range("N8").Value = 55.893224
$ws.Range("O8").Value = 0.1441466469015163
$ws.Range("P8").Value = 0.1441466469015162
$ws.Range("Q8").Value = 4.909573851144889
$ws.Range("R8").Value = 44.186164660304
$ws.Range("S8").Value = 0.0006877280704892797
$ws.Range("T8").Value = 0.0006877280704892796
$ws.Range("I9").Value = 0.00477103065019021
$ws.Range("J9").Value = 0.00477103065019021
$ws.Range("M9").Value = 24.41680433333333
$ws.Range("N9").Value = 73.25041299999999
$ws.Range("O9").Value = 0.1889102231444233
$ws.Range("P9").Value = 0.1889102231444233
$ws.Range("Q9").Value = 6.43420233283311
$ws.Range("R9").Value = 57.90782099549799
$ws.Range("S9").Value = 0.0009012964647563155
$ws.Range("T9").Value = 0.0009012964647563154
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.2809586666666666
$ws.Range("H10").Value = 0.842876
$ws.Range("I10").Value = 0.005086847862502274
$ws.Range("J10").Value = 0.005086847862502274
$ws.Range("M10").Value = 38.45264233333334
$ws.Range("N10").Value = 115.357927
$ws.Range("O10").Value = 0.2975040117664333
$ws.Range("P10").Value = 0.2975040117664332
$ws.Range("Q10").Value = 10.80360311978355
$ws.Range("R10").Value = 97.232428078052
$ws.Range("S10").Value = 0.001513357646339932
$ws.Range("T10").Value = 0.001513357646339932
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.2809586666666666
$ws.Range("H11").Value = 0.842876
$ws.Range("I11").Value = 0.005086847862502274
$ws.Range("J11").Value = 0.005086847862502274
$ws.Range("O11").Value = 0.3694391181876273
$ws.Range("P11").Value = 0.3694391181876272
$ws.Range("Q11").Value = 13.41586483531333
$ws.Range("R11").Value = 120.74278351782
$ws.Range("S11").Value = 0.001879280588677457
$ws.Range("T11").Value = 0.001879280588677456
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.2809586666666666
$ws.Range("H12").Value = 0.842876
$ws.Range("I12").Value = 0.005086847862502274
$ws.Range("J12").Value = 0.005086847862502274
$ws.Range("M12").Value = 18.63107466666667
$ws.Range("N12").Value = 55.893224
$ws.Range("O12").Value = 0.1441466469015163
$ws.Range("P12").Value = 0.1441466469015162
$ws.Range("Q12").Value = 5.234561896913777
$ws.Range("R12").Value = 47.111057072224
$ws.Range("S12").Value = 0.000733252062677848
$ws.Range("T12").Value = 0.0007332520626778479
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.2809586666666666
$ws.Range("H13").Value = 0.842876
$ws.Range("I13").Value = 0.005086847862502274
$ws.Range("J13").Value = 0.005086847862502274
$ws.Range("M13").Value = 24.41680433333333
$ws.Range("N13").Value = 73.25041299999999
$ws.Range("O13").Value = 0.1889102231444233
$ws.Range("P13").Value = 0.1889102231444233
$ws.Range("Q13").Value = 6.860112789754221
$ws.Range("R13").Value = 61.74101510778799
$ws.Range("S13").Value = 0.0009609575648070372
$ws.Range("T13").Value = 0.0009609575648070371
$ws.Range("G14").Value = 48.051317
$ws.Range("H14").Value = 144.153951
$ws.Range("I14").Value = 0.8699846923101473
$ws.Range("J14").Value = 0.8699846923101472
$ws.Range("M14").Value = 38.45264233333334
$ws.Range("N14").Value = 115.357927
$ws.Range("O14").Value = 0.2975040117664333
$ws.Range("P14").Value = 0.2975040117664332
$ws.Range("Q14").Value = 1847.70010624662
$ws.Range("R14").Value = 16629.30095621958
$ws.Range("S14").Value = 0.2588239361376549
$ws.Range("T14").Value = 0.2588239361376548
$ws.Range("G15").Value = 48.051317
$ws.Range("H15").Value = 144.153951
$ws.Range("I15").Value = 0.8699846923101473
$ws.Range("J15").Value = 0.8699846923101472
$ws.Range("O15").Value = 0.3694391181876273
$ws.Range("P15").Value = 0.3694391181876272
$ws.Range("Q15").Value = 2294.465522914855
$ws.Range("R15").Value = 20650.1897062337
$ws.Range("S15").Value = 0.321406377563795
$ws.Range("T15").Value = 0.321406377563795
$ws.Range("G16").Value = 48.051317
$ws.Range("H16").Value = 144.153951
$ws.Range("I16").Value = 0.8699846923101473
$ws.Range("J16").Value = 0.8699846923101472
$ws.Range("M16").Value = 18.63107466666667
$ws.Range("N16").Value = 55.893224
$ws.Range("O16").Value = 0.1441466469015163
$ws.Range("P16").Value = 0.1441466469015162
$ws.Range("Q16").Value = 895.2476748586694
$ws.Range("R16").Value = 8057.229073728025
$ws.Range("S16").Value = 0.1254053762521551
$ws.Range("T16").Value = 0.125405376252155
$ws.Range("G17").Value = 48.051317
$ws.Range("H17").Value = 144.153951
$ws.Range("I17").Value = 0.8699846923101473
$ws.Range("J17").Value = 0.8699846923101472
$ws.Range("M17").Value = 24.41680433333333
$ws.Range("N17").Value = 73.25041299999999
$ws.Range("O17").Value = 0.1889102231444233
$ws.Range("P17").Value = 0.1889102231444233
$ws.Range("Q17").Value = 1173.259605147974
$ws.Range("R17").Value = 10559.33644633176
$ws.Range("S17").Value = 0.1643490023565424
$ws.Range("T17").Value = 0.1643490023565423
